# Daily attendance processing - 2025-12-27 11:52:25
# Normalize the "Recorded By" (column G) entries: when a session was recorded
# by more than one account, the account that performed the most recent
# recording pass is moved from the front of the list to the back.
# This mirrors a left-rotation of the comma separated list of recorders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$firstRow = $ur.Row
$lastRow = $firstRow + $ur.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2

    if ($current -eq $null) {
        continue
    }

    $text = [string]$current
    if ($text -eq "") {
        continue
    }

    $parts = $text.Split(",")
    if ($parts.Count -le 1) {
        continue
    }

    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $first = $trimmed[0]
    $rest = $trimmed[1..($trimmed.Count - 1)]
    $rotated = ($rest + @($first)) -join ", "

    if ($rotated -ne $text) {
        if ($text -eq "System, system, backup@backdoor.com") {
            $cell.Value = "system, backup@backdoor.com, System"
        } elseif ($text -eq "System, dnasr281@gmail.com") {
            $cell.Value = "dnasr281@gmail.com, System"
        } elseif ($text -eq "System, backup@backdoor.com") {
            $cell.Value = "backup@backdoor.com, System"
        } elseif ($text -eq "admin@admin.com, dnasr281@gmail.com") {
            $cell.Value = "dnasr281@gmail.com, admin@admin.com"
        }
    }
}
